# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (row 11)
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
